# Update DNA extraction example values and ontology ref for sample description
$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": bump the template Version value ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.4"

# --- Sheet "Sample": fix header/ontology references for the sample description columns ---
$wsSample = $wb.Worksheets.Item("Sample")
$wsSample.Range("AF1").Value = "Characteristic [sample description]"
$wsSample.Range("AG1").Value = "Term Source REF (DPBO:0000175)"
$wsSample.Range("AH1").Value = "Term Accession Number (DPBO:0000175)"

# --- Keep the table column definitions (annotationTable) in sync with the header text ---
$table = $wsSample.ListObjects.Item("annotationTable")
$table.ListColumns.Item("Characteristic [Sample description]").Name = "Characteristic [sample description]"
$table.ListColumns.Item("Term Source REF (MIAPPE:0079)").Name = "Term Source REF (DPBO:0000175)"
$table.ListColumns.Item("Term Accession Number (MIAPPE:0079)").Name = "Term Accession Number (DPBO:0000175)"
